# Scheduled market-data refresh for the Leve-profit tracking workbook.
# Updates currentAveragePrice / NQ / HQ price & profit columns (H:N) across
# the eight crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with freshly pulled market-board figures. Row/column layout is untouched;
# only numeric leaf values change (plus two column(N) cells that now
# appear/disappear because the HQ-profit figure became computable/not).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 40
$ws.Range("H40").Value = 1987  # was 1987.25
$ws.Range("I40").Value = 1992.6666  # was 1993
$ws.Range("K40").Value = 1992.6666  # was 1993
$ws.Range("M40").Value = -1817.6666  # was -1818

$ws = $wb.Worksheets.Item("ALC")  # row 51
$ws.Range("H51").Value = 6636.364  # was 6842.7144
$ws.Range("I51").Value = 8000  # was 6800
$ws.Range("J51").Value = 6125  # was 6866.4443
$ws.Range("K51").Value = 8000  # was 6800
$ws.Range("L51").Value = 6125  # was 6866.4443
$ws.Range("M51").Value = -7516  # was -6316
$ws.Range("N51").Value = -7093  # was -7834.4443

$ws = $wb.Worksheets.Item("ALC")  # row 62
$ws.Range("H62").Value = 300  # was 250
$ws.Range("I62").Value = 300  # was 250
$ws.Range("K62").Value = 300  # was 250
$ws.Range("M62").Value = 324  # was 374

$ws = $wb.Worksheets.Item("ALC")  # row 65
$ws.Range("H65").Value = 300  # was 250
$ws.Range("I65").Value = 300  # was 250
$ws.Range("K65").Value = 1500  # was 1250
$ws.Range("M65").Value = 1620  # was 1870

$ws = $wb.Worksheets.Item("ALC")  # row 98
$ws.Range("H98").Value = 1109.2285  # was 1217.6857
$ws.Range("J98").Value = 1734.6666  # was 3000
$ws.Range("L98").Value = 1734.6666  # was 3000
$ws.Range("N98").Value = -4730.6666  # was -5996

$ws = $wb.Worksheets.Item("ALC")  # row 107
$ws.Range("H107").Value = 739.6667  # was 788.9286
$ws.Range("I107").Value = 392.23077  # was 420.75
$ws.Range("K107").Value = 392.23077  # was 420.75
$ws.Range("M107").Value = 1527.76923  # was 1499.25

$ws = $wb.Worksheets.Item("ALC")  # row 122
$ws.Range("H122").Value = 1109.2285  # was 1217.6857
$ws.Range("J122").Value = 1734.6666  # was 3000
$ws.Range("L122").Value = 5203.9998  # was 9000
$ws.Range("N122").Value = -10103.9998  # was -13900

$ws = $wb.Worksheets.Item("ALC")  # row 137
$ws.Range("H137").Value = 2009.6  # was 1835.8235
$ws.Range("I137").Value = 1916.5  # was 1694.4286
$ws.Range("J137").Value = 2071.6667  # was 1934.8
$ws.Range("K137").Value = 5749.5  # was 5083.2858
$ws.Range("L137").Value = 6215.000100000001  # was 5804.4
$ws.Range("M137").Value = -3199.5  # was -2533.2858
$ws.Range("N137").Value = -11315.0001  # was -10904.4

$ws = $wb.Worksheets.Item("ALC")  # row 138
$ws.Range("H138").Value = 3093.061  # was 3197.946
$ws.Range("I138").Value = 4093.65  # was 4138.65
$ws.Range("J138").Value = 2770.2903  # was 2849.537
$ws.Range("K138").Value = 12280.95  # was 12415.95
$ws.Range("L138").Value = 8310.8709  # was 8548.610999999999
$ws.Range("M138").Value = -7140.950000000001  # was -7275.949999999999
$ws.Range("N138").Value = -18590.8709  # was -18828.611

$ws = $wb.Worksheets.Item("ARM")  # row 2
$ws.Range("H2").Value = 214404.19  # was 199189.28
$ws.Range("I2").Value = 292900.6  # was 265138.2
$ws.Range("K2").Value = 292900.6  # was 265138.2
$ws.Range("M2").Value = -292787.6  # was -265025.2

$ws = $wb.Worksheets.Item("ARM")  # row 32
$ws.Range("H32").Value = 13471.295  # was 13583.742
$ws.Range("I32").Value = 9287.815000000001  # was 8822.366
$ws.Range("J32").Value = 20383.13  # was 21392.4
$ws.Range("K32").Value = 9287.815000000001  # was 8822.366
$ws.Range("L32").Value = 20383.13  # was 21392.4
$ws.Range("M32").Value = -9000.815000000001  # was -8535.366
$ws.Range("N32").Value = -20957.13  # was -21966.4

$ws = $wb.Worksheets.Item("ARM")  # row 74
$ws.Range("H74").Value = 789.44446  # was 803.2954999999999
$ws.Range("I74").Value = 587.6316  # was 598.6486
$ws.Range("K74").Value = 587.6316  # was 598.6486
$ws.Range("M74").Value = 286.3684  # was 275.3514

$ws = $wb.Worksheets.Item("ARM")  # row 77
$ws.Range("H77").Value = 789.44446  # was 803.2954999999999
$ws.Range("I77").Value = 587.6316  # was 598.6486
$ws.Range("K77").Value = 2938.158  # was 2993.243
$ws.Range("M77").Value = 1429.842  # was 1374.757

$ws = $wb.Worksheets.Item("ARM")  # row 82
$ws.Range("H82").Value = 100000  # was 75000
$ws.Range("J82").Value = 100000  # was 75000
$ws.Range("L82").Value = 100000  # was 75000
$ws.Range("N82").Value = -100722  # was -75722

$ws = $wb.Worksheets.Item("ARM")  # row 85
$ws.Range("H85").Value = 100000  # was 75000
$ws.Range("J85").Value = 100000  # was 75000
$ws.Range("L85").Value = 100000  # was 75000
$ws.Range("N85").Value = -102496  # was -77496

$ws = $wb.Worksheets.Item("ARM")  # row 112
$ws.Range("H112").Value = 49192  # was 48794.668
$ws.Range("J112").Value = 49192  # was 48794.668
$ws.Range("L112").Value = 49192  # was 48794.668
$ws.Range("N112").Value = -52146  # was -51748.668

$ws = $wb.Worksheets.Item("ARM")  # row 116
$ws.Range("H116").Value = 214404.19  # was 199189.28
$ws.Range("I116").Value = 292900.6  # was 265138.2
$ws.Range("K116").Value = 292900.6  # was 265138.2
$ws.Range("M116").Value = -290606.6  # was -262844.2

$ws = $wb.Worksheets.Item("ARM")  # row 124
$ws.Range("H124").Value = 13809.667  # was 15714.5
$ws.Range("J124").Value = 13809.667  # was 15714.5
$ws.Range("L124").Value = 13809.667  # was 15714.5
$ws.Range("N124").Value = -23629.667  # was -25534.5

$ws = $wb.Worksheets.Item("ARM")  # row 125
$ws.Range("H125").Value = 14750  # was 15000
$ws.Range("J125").Value = 14750  # was 15000
$ws.Range("L125").Value = 14750  # was 15000
$ws.Range("N125").Value = -24590  # was -24840

$ws = $wb.Worksheets.Item("ARM")  # row 132
$ws.Range("H132").Value = 2274.1667  # was 2380.7346
$ws.Range("I132").Value = 2082.48  # was 2248.318
$ws.Range("J132").Value = 2439.4138  # was 2488.6296
$ws.Range("K132").Value = 6247.440000000001  # was 6744.954000000001
$ws.Range("L132").Value = 7318.241399999999  # was 7465.888800000001
$ws.Range("M132").Value = -3717.440000000001  # was -4214.954000000001
$ws.Range("N132").Value = -12378.2414  # was -12525.8888

$ws = $wb.Worksheets.Item("BSM")  # row 3
$ws.Range("H3").Value = 214404.19  # was 199189.28
$ws.Range("I3").Value = 292900.6  # was 265138.2
$ws.Range("K3").Value = 292900.6  # was 265138.2
$ws.Range("M3").Value = -292786.6  # was -265024.2

$ws = $wb.Worksheets.Item("BSM")  # row 20
$ws.Range("H20").Value = 3319.8  # was 2262.375
$ws.Range("I20").Value = 2033.3334  # was 1266.6666
$ws.Range("K20").Value = 2033.3334  # was 1266.6666
$ws.Range("M20").Value = -1786.3334  # was -1019.6666

$ws = $wb.Worksheets.Item("BSM")  # row 99
$ws.Range("H99").Value = 744.25  # was 782
$ws.Range("J99").Value = 480  # was 0
$ws.Range("L99").Value = 480  # was 0
$ws.Range("N99").Value = -3476  # was None

$ws = $wb.Worksheets.Item("BSM")  # row 105
$ws.Range("H105").Value = 2306.2222  # was 2390.32
$ws.Range("I105").Value = 2111.2083  # was 2189.0454
$ws.Range("K105").Value = 2111.2083  # was 2189.0454
$ws.Range("M105").Value = -364.2082999999998  # was -442.0454

$ws = $wb.Worksheets.Item("BSM")  # row 134
$ws.Range("H134").Value = 7041.4346  # was 5915.8213
$ws.Range("I134").Value = 8642.866  # was 6396.8096
$ws.Range("J134").Value = 4038.75  # was 4472.857
$ws.Range("K134").Value = 25928.598  # was 19190.4288
$ws.Range("L134").Value = 12116.25  # was 13418.571
$ws.Range("M134").Value = -23393.598  # was -16655.4288
$ws.Range("N134").Value = -17186.25  # was -18488.571

$ws = $wb.Worksheets.Item("CRP")  # row 31
$ws.Range("H31").Value = 2589.7368  # was 2648.6216
$ws.Range("I31").Value = 1149.1364  # was 1184.2858
$ws.Range("K31").Value = 1149.1364  # was 1184.2858
$ws.Range("M31").Value = -854.1364000000001  # was -889.2858000000001

$ws = $wb.Worksheets.Item("CRP")  # row 34
$ws.Range("H34").Value = 2589.7368  # was 2648.6216
$ws.Range("I34").Value = 1149.1364  # was 1184.2858
$ws.Range("K34").Value = 1149.1364  # was 1184.2858
$ws.Range("M34").Value = -947.1364000000001  # was -982.2858000000001

$ws = $wb.Worksheets.Item("CRP")  # row 58
$ws.Range("H58").Value = 1062220.4  # was 1116642.2
$ws.Range("J58").Value = 1672.4117  # was 1714.8125
$ws.Range("L58").Value = 1672.4117  # was 1714.8125
$ws.Range("N58").Value = -2078.4117  # was -2120.8125

$ws = $wb.Worksheets.Item("CRP")  # row 132
$ws.Range("H132").Value = 2402.4285  # was 2548.7693
$ws.Range("I132").Value = 1237.3334  # was 1329.5
$ws.Range("K132").Value = 3712.0002  # was 3988.5
$ws.Range("M132").Value = -1182.0002  # was -1458.5

$ws = $wb.Worksheets.Item("CRP")  # row 134
$ws.Range("H134").Value = 1507.6333  # was 1285
$ws.Range("I134").Value = 1483.3158  # was 1173.0385
$ws.Range("K134").Value = 4449.9474  # was 3519.1155
$ws.Range("M134").Value = -1914.9474  # was -984.1155000000003

$ws = $wb.Worksheets.Item("CRP")  # row 136
$ws.Range("H136").Value = 1062220.4  # was 1116642.2
$ws.Range("J136").Value = 1672.4117  # was 1714.8125
$ws.Range("L136").Value = 5017.2351  # was 5144.4375
$ws.Range("N136").Value = -10117.2351  # was -10244.4375

$ws = $wb.Worksheets.Item("CRP")  # row 141
$ws.Range("H141").Value = 65266.383  # was 64532.645
$ws.Range("J141").Value = 64705.25  # was 63958.23
$ws.Range("L141").Value = 64705.25  # was 63958.23
$ws.Range("N141").Value = -75065.25  # was -74318.23000000001

$ws = $wb.Worksheets.Item("CUL")  # row 80
$ws.Range("H80").Value = 2495.3333  # was 2445.8
$ws.Range("J80").Value = 2495.3333  # was 2445.8
$ws.Range("L80").Value = 7485.999899999999  # was 7337.400000000001
$ws.Range("N80").Value = -9357.999899999999  # was -9209.400000000001

$ws = $wb.Worksheets.Item("CUL")  # row 83
$ws.Range("H83").Value = 2495.3333  # was 2445.8
$ws.Range("J83").Value = 2495.3333  # was 2445.8
$ws.Range("L83").Value = 22457.9997  # was 22012.2
$ws.Range("N83").Value = -31817.9997  # was -31372.2

$ws = $wb.Worksheets.Item("CUL")  # row 113
$ws.Range("H113").Value = 3962.2646  # was 3764.8918
$ws.Range("J113").Value = 981.44446  # was 1036.1
$ws.Range("L113").Value = 2944.33338  # was 3108.3
$ws.Range("N113").Value = -7284.33338  # was -7448.299999999999

$ws = $wb.Worksheets.Item("CUL")  # row 131
$ws.Range("H131").Value = 21129.674  # was 13083.974
$ws.Range("I131").Value = 422  # was 444
$ws.Range("J131").Value = 23101.834  # was 13425.595
$ws.Range("K131").Value = 1266  # was 1332
$ws.Range("L131").Value = 69305.50199999999  # was 40276.785
$ws.Range("M131").Value = 3774  # was 3708
$ws.Range("N131").Value = -79385.50199999999  # was -50356.785

$ws = $wb.Worksheets.Item("GSM")  # row 132
$ws.Range("H132").Value = 1136130.4  # was 1170521.6
$ws.Range("I132").Value = 1331284.8  # was 1484716.9
$ws.Range("J132").Value = 4234.8  # was 3510.4285
$ws.Range("K132").Value = 3993854.4  # was 4454150.699999999
$ws.Range("L132").Value = 12704.4  # was 10531.2855
$ws.Range("M132").Value = -3991324.4  # was -4451620.699999999
$ws.Range("N132").Value = -17764.4  # was -15591.2855

$ws = $wb.Worksheets.Item("GSM")  # row 135
$ws.Range("H135").Value = 65944.27  # was 78899.336
$ws.Range("J135").Value = 65944.27  # was 78899.336
$ws.Range("L135").Value = 65944.27  # was 78899.336
$ws.Range("N135").Value = -76084.27  # was -89039.336

$ws = $wb.Worksheets.Item("LTW")  # row 55
$ws.Range("H55").Value = 8334150  # was 7143571
$ws.Range("I55").Value = 33334000  # was 20000440
$ws.Range("K55").Value = 33334000  # was 20000440
$ws.Range("M55").Value = -33333827  # was -20000267

$ws = $wb.Worksheets.Item("LTW")  # row 122
$ws.Range("H122").Value = 4500  # was 5500
$ws.Range("I122").Value = 3900  # was 5000
$ws.Range("K122").Value = 11700  # was 15000
$ws.Range("M122").Value = -9250  # was -12550

$ws = $wb.Worksheets.Item("LTW")  # row 132
$ws.Range("H132").Value = 3810.682  # was 3378.8628
$ws.Range("I132").Value = 4142.1113  # was 3131.0386
$ws.Range("J132").Value = 3581.2307  # was 3636.6
$ws.Range("K132").Value = 12426.3339  # was 9393.1158
$ws.Range("L132").Value = 10743.6921  # was 10909.8
$ws.Range("M132").Value = -9896.333899999998  # was -6863.1158
$ws.Range("N132").Value = -15803.6921  # was -15969.8

$ws = $wb.Worksheets.Item("LTW")  # row 136
$ws.Range("H136").Value = 2001.1428  # was 1712.5
$ws.Range("I136").Value = 1833.8334  # was 1568.8889
$ws.Range("K136").Value = 5501.5002  # was 4706.6667
$ws.Range("M136").Value = -2951.5002  # was -2156.6667

$ws = $wb.Worksheets.Item("WVR")  # row 46
$ws.Range("H46").Value = 93966.664  # was 65499
$ws.Range("J46").Value = 93966.664  # was 65499
$ws.Range("L46").Value = 93966.664  # was 65499
$ws.Range("N46").Value = -94428.664  # was -65961

$ws = $wb.Worksheets.Item("WVR")  # row 62
$ws.Range("H62").Value = 3500  # was 4000
$ws.Range("I62").Value = 3500  # was 4000
$ws.Range("K62").Value = 3500  # was 4000
$ws.Range("M62").Value = -2876  # was -3376

$ws = $wb.Worksheets.Item("WVR")  # row 65
$ws.Range("H65").Value = 3500  # was 4000
$ws.Range("I65").Value = 3500  # was 4000
$ws.Range("K65").Value = 17500  # was 20000
$ws.Range("M65").Value = -14380  # was -16880

$ws = $wb.Worksheets.Item("WVR")  # row 119
$ws.Range("H119").Value = 0  # was 45000
$ws.Range("J119").Value = 0  # was 45000
$ws.Range("L119").Value = 0  # was 45000
$ws.Range("N119").ClearContents()  # was -54676

$ws = $wb.Worksheets.Item("WVR")  # row 126
$ws.Range("H126").Value = 3015.4285  # was 3408.2856
$ws.Range("I126").Value = 2977.4167  # was 3247.3845
$ws.Range("J126").Value = 3243.5  # was 5500
$ws.Range("K126").Value = 8932.250100000001  # was 9742.1535
$ws.Range("L126").Value = 9730.5  # was 16500
$ws.Range("M126").Value = -6462.250100000001  # was -7272.1535
$ws.Range("N126").Value = -14670.5  # was -21440

$ws = $wb.Worksheets.Item("WVR")  # row 132
$ws.Range("H132").Value = 1554.6154  # was 1613.8649
$ws.Range("I132").Value = 1383.6207  # was 1452.1482
$ws.Range("K132").Value = 4150.8621  # was 4356.444600000001
$ws.Range("M132").Value = -1620.8621  # was -1826.444600000001

$ws = $wb.Worksheets.Item("WVR")  # row 134
$ws.Range("H134").Value = 93966.664  # was 65499
$ws.Range("J134").Value = 93966.664  # was 65499
$ws.Range("L134").Value = 281899.992  # was 196497
$ws.Range("N134").Value = -286969.992  # was -201567

$ws = $wb.Worksheets.Item("WVR")  # row 136
$ws.Range("H136").Value = 13229199  # was 13551873
$ws.Range("I136").Value = 21368852  # was 23149574
$ws.Range("J136").Value = 2263.25  # was 2176.2942
$ws.Range("K136").Value = 64106556  # was 69448722
$ws.Range("L136").Value = 6789.75  # was 6528.882599999999
$ws.Range("M136").Value = -64104006  # was -69446172
$ws.Range("N136").Value = -11889.75  # was -11628.8826

